# "Update công nợ + phí ship"
#
# Sheet layout (1-based COM index -> name):
#   1 11.12.23   2 29.12.23   3 6.1.24   4 28.1.24
#   5 20.2.24    6 17.3.24    7 20.3.24
#
$wb = $excel.ActiveWorkbook

$ws4 = $wb.Worksheets.Item(4)   # 28.1.24
$ws5 = $wb.Worksheets.Item(5)   # 20.2.24
$ws6 = $wb.Worksheets.Item(6)   # 17.3.24
$ws7 = $wb.Worksheets.Item(7)   # 20.3.24

# ---------------------------------------------------------------------
# 20.2.24 (sheet 5): new công nợ line + note on the balance row
# ---------------------------------------------------------------------
$ws5.Range("A9").Value = "Chuyển khoản ngày 23.3.24"
$ws5.Range("B9").Value = 50000000
$ws5.Range("C9").Value = 1
$ws5.Range("D9").Formula = "=C9*B9"

$ws5.Range("E14").Value = "Chuyển sang ngày 17.3"

# Widen column E so the new note is readable.
$ws5.Columns.Item(5).ColumnWidth = 20.48

# ---------------------------------------------------------------------
# 17.3.24 (sheet 6): phí ship line + công nợ carried over from 20.2.24
# ---------------------------------------------------------------------
$ws6.Range("B9").Value = 400000
$ws6.Range("C9").Value = 1
$ws6.Range("D9").Formula = "=C9*B9"

$ws6.Range("A10").Value = "Chuyển công nợ ngày 20.2.24"
$ws6.Range("B10").Value = 49150000
$ws6.Range("C10").Value = 1
$ws6.Range("D10").Formula = "=C10*B10"

$ws6.Range("D11").Formula = "=SUM(D8:D10)"

# ---------------------------------------------------------------------
# 20.3.24 (sheet 7): Màn + Pin lines
# ---------------------------------------------------------------------
$ws7.Range("A8").Value = "Màn"
$ws7.Range("B8").Value = 200000
$ws7.Range("C8").Value = 1
$ws7.Range("D8").Formula = "=C8*B8"

$ws7.Range("A9").Value = "Pin "
$ws7.Range("B9").Value = 250000
$ws7.Range("C9").Value = 2
$ws7.Range("D9").Formula = "=C9*B9"

# ---------------------------------------------------------------------
# View-state: selections per sheet (last Activate()/Select() wins the
# workbook's active tab, so 20.2.24 must be selected last to become
# the active sheet, matching activeTab going from 20.3.24 to 20.2.24).
# ---------------------------------------------------------------------
$ws4.Range("A1:E16").Select()
$ws7.Range("E9").Select()
$ws6.Range("A1:E14").Select()
$ws5.Activate()
$ws5.Range("D21").Select()
